# TC06 - Buy Product
# Applies the changes described by the commit:
#  - updates the shared "email" value used on TC01 (and reused by the new sheet)
#  - adds a new worksheet TC06 at the end of the workbook with header/data rows
#  - makes TC06 the active sheet/tab (mirrors the previous "TC05 active" state)
#  - normalises a few cell selections on the existing sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the reused e-mail value (shared string) on TC01
# ---------------------------------------------------------------------------
$tc01 = $wb.Worksheets.Item("TC01")
$tc01.Range("B2").Value = "test0000@test.io"
# give this cell an explicit (no-op) style so it gets its own cell format,
# matching the distinct formatting applied to the e-mail cell in the source file
$tc01.Range("B2").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Create the new TC06 sheet ("Buy Product") at the end of the workbook
# ---------------------------------------------------------------------------
$tc05 = $wb.Worksheets.Item("TC05")
$tc06 = $wb.Worksheets.Add()
$tc06.Name = "TC06"

$tc06.Range("A1").Value = "baseUrl"
$tc06.Range("B1").Value = "userEmail"
$tc06.Range("C1").Value = "password"
$tc06.Range("D1").Value = "quantity"
$tc06.Range("E1").Value = "size"

$tc06.Range("A2").Value = "http://automationpractice.com/index.php"
$tc06.Range("B2").Value = "test0000@test.io"
$tc06.Range("C2").Value = "Pass1234"
$tc06.Range("D2").Value = 2
$tc06.Range("E2").Value = "L"

# match the formatting applied to the e-mail column (header + value)
$tc06.Range("B1").Style = "Normal"
$tc06.Range("B2").Style = "Normal"

$tc06.Columns.Item(1).ColumnWidth = 33.94
$tc06.Columns.Item(2).ColumnWidth = 15.33

# Move the new sheet to the end of the tab strip, after TC05
$tc06.Move($null, $tc05)

# ---------------------------------------------------------------------------
# 3. Normalise the selections on the previously-edited sheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("TC02").Range("A1").Select()
$wb.Worksheets.Item("TC03").Range("A2").Select()
$wb.Worksheets.Item("TC04").Range("A18").Select()
$wb.Worksheets.Item("TC05").Range("A1").Select()
$wb.Worksheets.Item("TC01").Range("B2").Select()

# ---------------------------------------------------------------------------
# 4. Activate TC06 (becomes the selected tab / active sheet, like TC05 was)
# ---------------------------------------------------------------------------
$tc06fresh = $wb.Worksheets.Item("TC06")
$tc06fresh.Range("E2").Select()
$tc06fresh.Activate()
